$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Lương")

# Remove the extra "Phụ cấp tại <location>" rows for LONG XUYÊN and SÓC TRĂNG.
# Row 13 = "Phụ cấp tại LONG XUYÊN", Row 24 = "Phụ cấp tại SÓC TRĂNG" (before any shifting).
# Delete row 24 first so row 13's index stays valid.
$ws.Rows.Item(24).Delete()
$ws.Rows.Item(13).Delete()
